# Regenerate the localization-status handback report:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on every sheet that surfaces it (Overview E2/F2, zh-cn C2, de-de C2).
#  - Latest Handback DateTime is refreshed for both locales.
#  - The de-de "handback out of date" Error Detail is cleared now that the
#    handback file is in sync (version-mismatch warning no longer applies).
#  - The Status / Error Detail columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-20 10:54:13"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-20 10:54:19"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (Status / Error Detail columns widened) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527       # C: Status
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839      # P: Error Detail

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527       # C: Status
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839      # P: Error Detail
